$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header
$ws.Range("A1").Value = "Cluster"

# Update data rows 2-27 with the refreshed cluster list/values
$ws.Range("A2").Value = '3153 Sacred Heart Community St Kilda Tier 1A'
$ws.Range("B2").Value = 21
$ws.Range("A3").Value = '3600 Belvedere Aged Care Noble Park'
$ws.Range("B3").Value = 39
$ws.Range("A4").Value = '3612 BlueCross Glengowrie'
$ws.Range("B4").Value = 56
$ws.Range("A5").Value = '3684 Homestyle Aged Care Langford Grange Cranbourne East'
$ws.Range("B5").Value = 30
$ws.Range("A6").Value = '3824 Estia Health South Morang'
$ws.Range("B6").Value = 32
$ws.Range("A7").Value = '3980 Arcare Keysborough Aged Care Keysborough'
$ws.Range("B7").Value = 27
$ws.Range("A8").Value = '4075 Homestyle Aged Care Ferndale Gardens Aged Care Services Bayswater North'
$ws.Range("B8").Value = 16
$ws.Range("A9").Value = '4518 Regis Aged Care Fawkner'
$ws.Range("B9").Value = 21
$ws.Range("A10").Value = 'ACFS Port Logistics Altona'
$ws.Range("B10").Value = 14
$ws.Range("A11").Value = 'Bespoke Childcare Dingley Village'
$ws.Range("B11").Value = 13
$ws.Range("A12").Value = 'Community Kids Pascoe Vale Early Education Centre Pascoe Vale'
$ws.Range("B12").Value = 23
$ws.Range("A13").Value = 'Guardian Childcare Caulfield'
$ws.Range("B13").Value = 17
$ws.Range("A14").Value = 'Hello Fresh Warehouse Ravenhall'
$ws.Range("B14").Value = 18
$ws.Range("A15").Value = 'Honeyeater Hairdressers Bendigo'
$ws.Range("B15").Value = 10
$ws.Range("A16").Value = 'Inghams Enterprises Somerville'
$ws.Range("B16").Value = 27
$ws.Range("A17").Value = 'JBS Australia Brooklyn'
$ws.Range("B17").Value = 10
$ws.Range("A18").Value = 'Metro Tunnel Shine Domain Site Albert Road Construction Site South Melbourne'
$ws.Range("B18").Value = 10
$ws.Range("A19").Value = 'Northern Health The Northern Hospital Epping'
$ws.Range("B19").Value = 17
$ws.Range("A20").Value = 'Shawlands Caravan Park Dandenong South'
$ws.Range("B20").Value = 17
$ws.Range("A21").Value = 'St Vincents Hospital Emergency Department Melbourne'
$ws.Range("B21").Value = 32
$ws.Range("A22").Value = 'The Robin Hood Inn Drouin West'
$ws.Range("B22").Value = 49
$ws.Range("A23").Value = 'The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1A'
$ws.Range("B23").Value = 10
$ws.Range("A24").Value = 'The Royal Children''s Hospital Parkville'
$ws.Range("B24").Value = 11
$ws.Range("A25").Value = 'Visy Recycling Springvale'
$ws.Range("B25").Value = 10
$ws.Range("A26").Value = 'Werribee Mercy Hospital Emergency Department'
$ws.Range("B26").Value = 33
$ws.Range("A27").Value = 'Western Health Sunshine Hospital Emergency Department'
$ws.Range("B27").Value = 22

# Delete now-unused rows 28-29 (table shrank from 29 to 27 rows), so the
# used range/dimension shrinks to B27
$ws.Range("28:29").Delete()
